# Huge CSS Revamp, added first test checking all grids
#
# Content change: the "Adres aplikacji" hyperlink cell (B38) used to show the
# internal test URL; it now shows the public T-Mobile PL URL instead.
# (The underlying hyperlink target itself is left untouched - only the
# displayed text/shared-string changes, matching the source diff which only
# touches xl/sharedStrings.xml + the cached <v> index in the sheet, not the
# worksheet's hyperlink relationships file.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B38").Value = "http://www.t-mobile.pl/"

# View state: the workbook was left scrolled down with B38 selected.
$ws.Range("B38").Select()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
